$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 324 (existing data shifts down by 2).
$ws.Rows("324:325").Insert()

# New row 324 ("Primera" quality) - latest weekly price entry.
$ws.Cells.Item(324, 1).Value = 11
$ws.Cells.Item(324, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(324, 3).Value = "Bíobío"
$ws.Cells.Item(324, 4).Value = 45027
$ws.Cells.Item(324, 5).Value = 8
$ws.Cells.Item(324, 6).Value = 100112017
$ws.Cells.Item(324, 7).Value = "Apio"
$ws.Cells.Item(324, 8).Value = "Americana (o)"
$ws.Cells.Item(324, 9).Value = "Primera"
$ws.Cells.Item(324, 10).Value = 350
$ws.Cells.Item(324, 11).Value = 6500
$ws.Cells.Item(324, 12).Value = 7000
$ws.Cells.Item(324, 13).Value = 6714
$ws.Cells.Item(324, 14).Value = "$/docena de matas"
$ws.Cells.Item(324, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(324, 16).Value = 1119
$ws.Cells.Item(324, 17).Value = 6
$ws.Cells.Item(324, 18).Value = "Hortaliza"

# New row 325 ("Segunda" quality) - latest weekly price entry.
$ws.Cells.Item(325, 1).Value = 11
$ws.Cells.Item(325, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(325, 3).Value = "Bíobío"
$ws.Cells.Item(325, 4).Value = 45027
$ws.Cells.Item(325, 5).Value = 8
$ws.Cells.Item(325, 6).Value = 100112017
$ws.Cells.Item(325, 7).Value = "Apio"
$ws.Cells.Item(325, 8).Value = "Americana (o)"
$ws.Cells.Item(325, 9).Value = "Segunda"
$ws.Cells.Item(325, 10).Value = 200
$ws.Cells.Item(325, 11).Value = 5500
$ws.Cells.Item(325, 12).Value = 5500
$ws.Cells.Item(325, 13).Value = 5500
$ws.Cells.Item(325, 14).Value = "$/docena de matas"
$ws.Cells.Item(325, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(325, 16).Value = 917
$ws.Cells.Item(325, 17).Value = 6
$ws.Cells.Item(325, 18).Value = "Hortaliza"
